# update new orleans xlsx files
#
# 1. Insert a new "State" column into hotel_info (right after Hotel_Name,
#    before City) and populate it with "Louisiana" for the existing hotel row.
# 2. Re-order the worksheet tabs so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

$wsHotel  = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# --- 1. Add the new State column to hotel_info -----------------------------
# Hotel_Name is column B, City is column C -> insert a new column at C.
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"

# --- 2. Put review_info ahead of hotel_info in the tab order ---------------
$wsReview.Move($wsHotel)
